# Visit List Test Plan work: add Encounter columns (H:S) to the "Chart" sheet,
# style the new headers, format the new date column, and switch the active
# tab/selection over to the Chart sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# --- New header cells (row 1) ---------------------------------------------
# H1:K1 and P1:S1 get the "new" bold+white-fill header style.
$hdrNewRange = @("H1","I1","J1","K1")
$hdrNewValues = @("Enc1_Time","Enc1_AttPhy","Enc1_Type","Enc1_Reason")
for ($i = 0; $i -lt $hdrNewRange.Length; $i++) {
    $c = $ws.Range($hdrNewRange[$i])
    $c.Value = $hdrNewValues[$i]
    $c.Font.Bold = $true
    $c.Interior.ThemeColor = 2
}

# L1:O1 reuse the existing bold/green header style (same as A1:G1), so simply
# assigning the value (inheriting the row's default format) is enough.
$ws.Range("L1").Value = "Enc2_Time"
$ws.Range("M1").Value = "Enc2_AttPhy"
$ws.Range("N1").Value = "Enc2_Type"
$ws.Range("O1").Value = "Enc2_Reason"

$hdrNewRange2 = @("P1","Q1","R1","S1")
$hdrNewValues2 = @("Enc3_Time","Enc3_AttPhy","Enc3_Type","Enc3_Reason")
for ($i = 0; $i -lt $hdrNewRange2.Length; $i++) {
    $c = $ws.Range($hdrNewRange2[$i])
    $c.Value = $hdrNewValues2[$i]
    $c.Font.Bold = $true
    $c.Interior.ThemeColor = 2
}

# --- New data cell (row 2) --------------------------------------------------
$ws.Range("H2").NumberFormat = "m/d/yy h:mm"
$ws.Range("H2").Value = "Chest pain"

# --- Column widths for the new columns -------------------------------------
$ws.Columns.Item(8).ColumnWidth = 13 + 1/12      # H  -> 14
$ws.Columns.Item(9).ColumnWidth = 14.75          # I  -> ~15.6640625
$ws.Columns.Item(11).ColumnWidth = 11 + 7/12     # K  -> 12.5
$ws.Columns.Item(13).ColumnWidth = 12 + 7/12     # M  -> 13.5
$ws.Columns.Item(14).ColumnWidth = 11 + 11/12    # N  -> ~12.83203125
$ws.Columns.Item(15).ColumnWidth = 11 + 11/12    # O  -> ~12.83203125
$ws.Columns.Item(16).ColumnWidth = 16.75         # P  -> ~17.6640625
$ws.Columns.Item(17).ColumnWidth = 14 + 5/12     # Q  -> ~15.33203125
$ws.Columns.Item(19).ColumnWidth = 12 + 1/12     # S  -> 13

# --- Switch the active sheet/selection to Chart -----------------------------
$ws.Activate()
$ws.Range("H2").Select()
